$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 93: Course, Hours, Notes for the "diamond and hollow diamond" entry
$ws.Range("B93").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C93").Value = 0.75
$ws.Range("D93").Value = "Finish 1 small problem"

# Move the active selection to C93 (matches the author's final cursor position)
[void]$ws.Range("C93").Select()
